$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "68.437.56"
$ws.Range("E2").Value = "  +1.04%  "
Set-TextValue "D3" "3.770.32"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "595.73"
$ws.Range("E5").Value = "  -0.29%  "
Set-TextValue "D6" "168.61"
$ws.Range("E6").Value = "  -0.39%  "
Set-TextValue "D7" "3.769.61"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("E12").Value = "  -2.55%  "
Set-TextValue "D13" "0.0000264"
$ws.Range("E13").Value = "  -3.53%  "
Set-TextValue "D14" "36.57"
Set-TextValue "D15" "4.406.26"
$ws.Range("E15").Value = "  -0.33%  "
Set-TextValue "D16" "3.780.56"
$ws.Range("E16").Value = "  +0.00%  "
Set-TextValue "D17" "68.432.53"
$ws.Range("E17").Value = "  +1.00%  "
Set-TextValue "D18" "18.23"
$ws.Range("E18").Value = "  -4.15%  "
Set-TextValue "D19" "7.06"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("E20").Value = "  -0.41%  "
Set-TextValue "D21" "10.79"
$ws.Range("E21").Value = "  +2.11%  "
Set-TextValue "D22" "469.41"
$ws.Range("E22").Value = "  +0.74%  "
Set-TextValue "D23" "0.702"
$ws.Range("E23").Value = "  -3.38%  "
Set-TextValue "D24" "84.31"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("E26").Value = "  -0.15%  "
Set-TextValue "D27" "12.24"
$ws.Range("E27").Value = "  +0.36%  "
Set-TextValue "D28" "10.23"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  +0.12%  "
Set-TextValue "D30" "3.918.24"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  -4.34%  "
Set-TextValue "D32" "7.44"
$ws.Range("E32").Value = "  -2.50%  "
Set-TextValue "D33" "2.22"
$ws.Range("E33").Value = "  -1.56%  "
Set-TextValue "D34" "30.09"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("E35").Value = "  +0.91%  "
Set-TextValue "D36" "0.999"
Set-TextValue "D37" "3.726.35"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("E38").Value = "  -3.41%  "
Set-TextValue "D39" "3.48"
$ws.Range("E39").Value = "  -9.27%  "
Set-TextValue "D40" "0.138"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -0.33%  "
Set-TextValue "D42" "5.82"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -3.03%  "
Set-TextValue "D46" "1.97"
$ws.Range("E46").Value = "  +0.10%  "
Set-TextValue "D47" "43.55"
$ws.Range("E47").Value = "  +12.35%  "
$ws.Range("E48").Value = "  -1.64%  "
Set-TextValue "D49" "406.83"
$ws.Range("E49").Value = "  -0.13%  "
Set-TextValue "D50" "45.37"
$ws.Range("E50").Value = "  -2.00%  "
Set-TextValue "D51" "144.92"
$ws.Range("E51").Value = "  +1.97%  "
